## Apply the "Points Calculator" sheet addition to the workbook.
$wb = $excel.ActiveWorkbook

# --- 1. Insert the new worksheet right after "Tmin Tmax" ------------------
$afterSheet = $wb.Worksheets.Item("Tmin Tmax")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = "Points Calculator"

# --- 2. Header row (row 1) -------------------------------------------------
$headers = @(
    "sweep",
    "gear_ratio",
    "points",
    "endurance_battery_capacity",
    "endurance time",
    "endurance total time",
    "autocross time",
    "skidpad time",
    "acceleration time",
    "endurance points",
    "autocross points",
    "skidpad points",
    "acceleration points",
    "Total Points"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# Whole sheet defaults to the 11pt "Normal 2" look used throughout this tab.
$ws.Cells.Font.Size = 11

# "Total Points" column is bold, both header and the formula beneath it.
$ws.Range("N1").Font.Bold = $true
$ws.Range("N2").Font.Bold = $true

# --- 3. Formulas for row 2 ---------------------------------------------------
$ws.Range("J2").Formula = "=250*((('Tmin Tmax'!`$J`$11/F2)-1)/(('Tmin Tmax'!`$J`$11/'Tmin Tmax'!`$J`$5)-1))"
$ws.Range("K2").Formula = "=118.5*((('Tmin Tmax'!`$J`$10/G2)-1)/(('Tmin Tmax'!`$J`$10/'Tmin Tmax'!`$J`$4)-1))+6.5"
$ws.Range("L2").Formula = "=71.5*(((('Tmin Tmax'!`$J`$9/H2)^2)-1)/((('Tmin Tmax'!`$J`$9/'Tmin Tmax'!`$J`$3)^2)-1))+3.5"
$ws.Range("M2").Formula = "=95.5*((('Tmin Tmax'!`$J`$8/I2)-1)/(('Tmin Tmax'!`$J`$8/'Tmin Tmax'!`$J`$2)-1))+4.5"
$ws.Range("N2").Formula = "=SUM(J2:M2)"

# --- 4. Column widths / hidden columns --------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.8
$ws.Columns.Item(2).ColumnWidth = 9.4
$ws.Range("C1:D1").EntireColumn.Hidden = $true
$ws.Columns.Item(5).ColumnWidth = 13.6
$ws.Columns.Item(6).ColumnWidth = 16.4
$ws.Columns.Item(7).ColumnWidth = 12.9
$ws.Columns.Item(8).ColumnWidth = 11.1
$ws.Columns.Item(9).ColumnWidth = 14.8
$ws.Columns.Item(10).ColumnWidth = 14.9
$ws.Columns.Item(11).ColumnWidth = 14.1
$ws.Columns.Item(12).ColumnWidth = 12.5
$ws.Columns.Item(13).ColumnWidth = 16.2
$ws.Columns.Item(14).ColumnWidth = 9.8

# --- 5. Number formats -------------------------------------------------------
# gear_ratio column -> 2 decimal places
$ws.Range("B2:B15").NumberFormat = "0.00"
# the five time columns -> 3 decimal places
$ws.Range("E2:I15").NumberFormat = "0.000"
# Total Points -> 3 decimal places too
$ws.Range("N1:N2").NumberFormat = "0.000"

$ws.Range("K2").Select()

# --- 6. Activate this sheet so the tab selection matches -------------------
$ws.Activate()
